# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# market-data runner. Values only; no formulas/formatting involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 188.11111
$ws.Cells.Item(6, 9).Value = 207.16667
$ws.Cells.Item(6, 11).Value = 621.50001
$ws.Cells.Item(6, 13).Value = -509.50001
$ws.Cells.Item(19, 8).Value = 2289.1333
$ws.Cells.Item(19, 9).Value = 2063.4285
$ws.Cells.Item(19, 10).Value = 2486.625
$ws.Cells.Item(19, 11).Value = 2063.4285
$ws.Cells.Item(19, 12).Value = 2486.625
$ws.Cells.Item(19, 13).Value = -1888.4285
$ws.Cells.Item(19, 14).Value = -2836.625
$ws.Cells.Item(74, 8).Value = 2862
$ws.Cells.Item(74, 9).Value = 2862
$ws.Cells.Item(74, 11).Value = 2862
$ws.Cells.Item(74, 13).Value = -1926
$ws.Cells.Item(77, 8).Value = 2862
$ws.Cells.Item(77, 9).Value = 2862
$ws.Cells.Item(77, 11).Value = 14310
$ws.Cells.Item(77, 13).Value = -9630
$ws.Cells.Item(100, 8).Value = 2998.75
$ws.Cells.Item(100, 9).Value = 3331.6667
$ws.Cells.Item(100, 10).Value = 2000
$ws.Cells.Item(100, 11).Value = 3331.6667
$ws.Cells.Item(100, 12).Value = 2000
$ws.Cells.Item(100, 13).Value = -2790.6667
$ws.Cells.Item(100, 14).Value = -3082
$ws.Cells.Item(137, 8).Value = 3992.8
$ws.Cells.Item(137, 9).Value = 3632.5715
$ws.Cells.Item(137, 11).Value = 10897.7145
$ws.Cells.Item(137, 13).Value = -8347.7145
$ws.Cells.Item(138, 8).Value = 4213.8887
$ws.Cells.Item(138, 10).Value = 4908.483
$ws.Cells.Item(138, 12).Value = 14725.449
$ws.Cells.Item(138, 14).Value = -25005.449

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3344.4
$ws.Cells.Item(2, 9).Value = 3200.4783
$ws.Cells.Item(2, 10).Value = 4999.5
$ws.Cells.Item(2, 11).Value = 3200.4783
$ws.Cells.Item(2, 12).Value = 4999.5
$ws.Cells.Item(2, 13).Value = -3087.4783
$ws.Cells.Item(2, 14).Value = -5225.5
$ws.Cells.Item(74, 8).Value = 2827.7144
$ws.Cells.Item(74, 9).Value = 2966.3333
$ws.Cells.Item(74, 10).Value = 2723.75
$ws.Cells.Item(74, 11).Value = 2966.3333
$ws.Cells.Item(74, 12).Value = 2723.75
$ws.Cells.Item(74, 13).Value = -2092.3333
$ws.Cells.Item(74, 14).Value = -4471.75
$ws.Cells.Item(77, 8).Value = 2827.7144
$ws.Cells.Item(77, 9).Value = 2966.3333
$ws.Cells.Item(77, 10).Value = 2723.75
$ws.Cells.Item(77, 11).Value = 14831.6665
$ws.Cells.Item(77, 12).Value = 13618.75
$ws.Cells.Item(77, 13).Value = -10463.6665
$ws.Cells.Item(77, 14).Value = -22354.75
$ws.Cells.Item(110, 8).Value = 2159.375
$ws.Cells.Item(110, 9).Value = 2182.1428
$ws.Cells.Item(110, 11).Value = 2182.1428
$ws.Cells.Item(110, 13).Value = -137.1428000000001
$ws.Cells.Item(116, 8).Value = 3344.4
$ws.Cells.Item(116, 9).Value = 3200.4783
$ws.Cells.Item(116, 10).Value = 4999.5
$ws.Cells.Item(116, 11).Value = 3200.4783
$ws.Cells.Item(116, 12).Value = 4999.5
$ws.Cells.Item(116, 13).Value = -906.4783000000002
$ws.Cells.Item(116, 14).Value = -9587.5
$ws.Cells.Item(124, 8).Value = 27771.4
$ws.Cells.Item(124, 10).Value = 27771.4
$ws.Cells.Item(124, 12).Value = 27771.4
$ws.Cells.Item(124, 14).Value = -37591.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3344.4
$ws.Cells.Item(3, 9).Value = 3200.4783
$ws.Cells.Item(3, 10).Value = 4999.5
$ws.Cells.Item(3, 11).Value = 3200.4783
$ws.Cells.Item(3, 12).Value = 4999.5
$ws.Cells.Item(3, 13).Value = -3086.4783
$ws.Cells.Item(3, 14).Value = -5227.5
$ws.Cells.Item(86, 8).Value = 2085.1538
$ws.Cells.Item(86, 9).Value = 1741.2
$ws.Cells.Item(86, 10).Value = 3231.6667
$ws.Cells.Item(86, 11).Value = 1741.2
$ws.Cells.Item(86, 12).Value = 3231.6667
$ws.Cells.Item(86, 13).Value = -618.2
$ws.Cells.Item(86, 14).Value = -5477.6667
$ws.Cells.Item(89, 8).Value = 2085.1538
$ws.Cells.Item(89, 9).Value = 1741.2
$ws.Cells.Item(89, 10).Value = 3231.6667
$ws.Cells.Item(89, 11).Value = 8706
$ws.Cells.Item(89, 12).Value = 16158.3335
$ws.Cells.Item(89, 13).Value = -3090
$ws.Cells.Item(89, 14).Value = -27390.3335
$ws.Cells.Item(99, 8).Value = 4156.2915
$ws.Cells.Item(99, 9).Value = 4034.158
$ws.Cells.Item(99, 11).Value = 4034.158
$ws.Cells.Item(99, 13).Value = -2536.158
$ws.Cells.Item(105, 8).Value = 2520.5715
$ws.Cells.Item(105, 9).Value = 2521.8
$ws.Cells.Item(105, 11).Value = 2521.8
$ws.Cells.Item(105, 13).Value = -774.8000000000002
$ws.Cells.Item(107, 8).Value = 1023.6
$ws.Cells.Item(107, 9).Value = 1023.6
$ws.Cells.Item(107, 11).Value = 1023.6
$ws.Cells.Item(107, 13).Value = 896.4
$ws.Cells.Item(134, 8).Value = 2320.1538
$ws.Cells.Item(134, 9).Value = 2268.6
$ws.Cells.Item(134, 11).Value = 6805.799999999999
$ws.Cells.Item(134, 13).Value = -4270.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 12244.632
$ws.Cells.Item(31, 9).Value = 12636.25
$ws.Cells.Item(31, 11).Value = 12636.25
$ws.Cells.Item(31, 13).Value = -12341.25
$ws.Cells.Item(34, 8).Value = 12244.632
$ws.Cells.Item(34, 9).Value = 12636.25
$ws.Cells.Item(34, 11).Value = 12636.25
$ws.Cells.Item(34, 13).Value = -12434.25
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 30763.5
$ws.Cells.Item(81, 10).Value = 30763.5
$ws.Cells.Item(81, 12).Value = 30763.5
$ws.Cells.Item(81, 14).Value = -32759.5
$ws.Cells.Item(84, 8).Value = 30763.5
$ws.Cells.Item(84, 10).Value = 30763.5
$ws.Cells.Item(84, 12).Value = 92290.5
$ws.Cells.Item(84, 14).Value = -102274.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 384.2857
$ws.Cells.Item(26, 9).Value = 140
$ws.Cells.Item(26, 10).Value = 425
$ws.Cells.Item(26, 11).Value = 420
$ws.Cells.Item(26, 12).Value = 1275
$ws.Cells.Item(26, 13).Value = -132
$ws.Cells.Item(26, 14).Value = -1851
$ws.Cells.Item(68, 8).Value = 889
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 889
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(88, 8).Value = 999
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 8).Value = 999
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 516.375
$ws.Cells.Item(107, 9).Value = 196.25
$ws.Cells.Item(107, 11).Value = 588.75
$ws.Cells.Item(107, 13).Value = 1331.25
$ws.Cells.Item(122, 8).Value = 22766.666
$ws.Cells.Item(122, 10).Value = 17400
$ws.Cells.Item(122, 12).Value = 156600
$ws.Cells.Item(122, 14).Value = -161500
$ws.Cells.Item(140, 8).Value = 1978.2051
$ws.Cells.Item(140, 9).Value = 1776.8529
$ws.Cells.Item(140, 10).Value = 3347.4
$ws.Cells.Item(140, 11).Value = 5330.5587
$ws.Cells.Item(140, 12).Value = 10042.2
$ws.Cells.Item(140, 13).Value = -150.5587000000005
$ws.Cells.Item(140, 14).Value = -20402.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2975.8647
$ws.Cells.Item(132, 9).Value = 2840.6333
$ws.Cells.Item(132, 11).Value = 8521.8999
$ws.Cells.Item(132, 13).Value = -5991.8999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 22201.5
$ws.Cells.Item(7, 10).Value = 39748.25
$ws.Cells.Item(7, 12).Value = 39748.25
$ws.Cells.Item(7, 14).Value = -39972.25
$ws.Cells.Item(46, 8).Value = 7432.8
$ws.Cells.Item(46, 9).Value = 7999.222
$ws.Cells.Item(46, 10).Value = 6583.1665
$ws.Cells.Item(46, 11).Value = 7999.222
$ws.Cells.Item(46, 12).Value = 6583.1665
$ws.Cells.Item(46, 13).Value = -7811.222
$ws.Cells.Item(46, 14).Value = -6959.1665
$ws.Cells.Item(122, 8).Value = 5832.9165
$ws.Cells.Item(122, 9).Value = 5761.3335
$ws.Cells.Item(122, 11).Value = 17284.0005
$ws.Cells.Item(122, 13).Value = -14834.0005
$ws.Cells.Item(126, 8).Value = 22201.5
$ws.Cells.Item(126, 10).Value = 39748.25
$ws.Cells.Item(126, 12).Value = 119244.75
$ws.Cells.Item(126, 14).Value = -124184.75
$ws.Cells.Item(127, 8).Value = 66404.5
$ws.Cells.Item(127, 10).Value = 66404.5
$ws.Cells.Item(127, 12).Value = 66404.5
$ws.Cells.Item(127, 14).Value = -76324.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2284
$ws.Cells.Item(81, 9).Value = 2462.1304
$ws.Cells.Item(81, 10).Value = 1698.7142
$ws.Cells.Item(81, 11).Value = 4924.2608
$ws.Cells.Item(81, 12).Value = 3397.4284
$ws.Cells.Item(81, 13).Value = -3863.2608
$ws.Cells.Item(81, 14).Value = -5519.4284
$ws.Cells.Item(84, 8).Value = 2284
$ws.Cells.Item(84, 9).Value = 2462.1304
$ws.Cells.Item(84, 10).Value = 1698.7142
$ws.Cells.Item(84, 11).Value = 24621.304
$ws.Cells.Item(84, 12).Value = 16987.142
$ws.Cells.Item(84, 13).Value = -19317.304
$ws.Cells.Item(84, 14).Value = -27595.142
$ws.Cells.Item(107, 8).Value = 1921.75
$ws.Cells.Item(107, 9).Value = 2316.3333
$ws.Cells.Item(107, 10).Value = 1752.6428
$ws.Cells.Item(107, 11).Value = 6948.999899999999
$ws.Cells.Item(107, 12).Value = 5257.928400000001
$ws.Cells.Item(107, 13).Value = -5028.999899999999
$ws.Cells.Item(107, 14).Value = -9097.928400000001
$ws.Cells.Item(132, 8).Value = 3811.8096
$ws.Cells.Item(132, 9).Value = 3811.8096
$ws.Cells.Item(132, 11).Value = 11435.4288
$ws.Cells.Item(132, 13).Value = -8905.4288
